# Auto-generated edit script applying numeric updates to the Famfrit_Profits workbook
# Mirrors a scheduled market-data refresh: currentAveragePrice* / Leve price & profit columns
# are recomputed per leve row across all eight job sheets.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1449.8889
$ws.Range("I98").Value = 1449.8889
$ws.Range("K98").Value = 1449.8889
$ws.Range("M98").Value = 48.11110000000008
$ws.Range("H118").Value = 472.3
$ws.Range("I118").Value = 517.875
$ws.Range("J118").Value = 290
$ws.Range("K118").Value = 1553.625
$ws.Range("L118").Value = 870
$ws.Range("M118").Value = 103.375
$ws.Range("N118").Value = -4184
$ws.Range("H122").Value = 1449.8889
$ws.Range("I122").Value = 1449.8889
$ws.Range("K122").Value = 4349.6667
$ws.Range("M122").Value = -1899.6667
$ws.Range("H127").Value = 10240.75
$ws.Range("I127").Value = 469.5
$ws.Range("J127").Value = 13497.833
$ws.Range("K127").Value = 1408.5
$ws.Range("L127").Value = 40493.499
$ws.Range("M127").Value = 3551.5
$ws.Range("N127").Value = -50413.499
$ws.Range("H137").Value = 19701.13
$ws.Range("I137").Value = 23499.4
$ws.Range("J137").Value = 12579.375
$ws.Range("K137").Value = 70498.20000000001
$ws.Range("L137").Value = 37738.125
$ws.Range("M137").Value = -67948.20000000001
$ws.Range("N137").Value = -42838.125

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 23811732
$ws.Range("I122").Value = 1592.1
$ws.Range("J122").Value = 83337090
$ws.Range("K122").Value = 4776.299999999999
$ws.Range("L122").Value = 250011270
$ws.Range("M122").Value = -2326.299999999999
$ws.Range("N122").Value = -250016170

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 5556.4443
$ws.Range("I99").Value = 3500
$ws.Range("K99").Value = 3500
$ws.Range("M99").Value = -2002
$ws.Range("H107").Value = 2345.3333
$ws.Range("I107").Value = 2345.3333
$ws.Range("K107").Value = 2345.3333
$ws.Range("M107").Value = -425.3332999999998
$ws.Range("H128").Value = 3000
$ws.Range("I128").Value = 3000
$ws.Range("K128").Value = 9000
$ws.Range("M128").Value = -6510

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 6398000
$ws.Range("I4").Value = 6398000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 6398000
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -6397888
$ws.Range("N4").ClearContents()
$ws.Range("H10").Value = 400537.4
$ws.Range("I10").Value = 500421.75
$ws.Range("K10").Value = 500421.75
$ws.Range("M10").Value = -500282.75
$ws.Range("H31").Value = 7815255.5
$ws.Range("I31").Value = 2024.24
$ws.Range("J31").Value = 35719652
$ws.Range("K31").Value = 2024.24
$ws.Range("L31").Value = 35719652
$ws.Range("M31").Value = -1729.24
$ws.Range("N31").Value = -35720242
$ws.Range("H34").Value = 7815255.5
$ws.Range("I34").Value = 2024.24
$ws.Range("J34").Value = 35719652
$ws.Range("K34").Value = 2024.24
$ws.Range("L34").Value = 35719652
$ws.Range("M34").Value = -1822.24
$ws.Range("N34").Value = -35720056
$ws.Range("H132").Value = 4240.55
$ws.Range("I132").Value = 4178.5
$ws.Range("K132").Value = 12535.5
$ws.Range("M132").Value = -10005.5
$ws.Range("H134").Value = 3094.5
$ws.Range("I134").Value = 2707.2
$ws.Range("J134").Value = 4385.5
$ws.Range("K134").Value = 8121.599999999999
$ws.Range("L134").Value = 13156.5
$ws.Range("M134").Value = -5586.599999999999
$ws.Range("N134").Value = -18226.5

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 1136.4736
$ws.Range("I39").Value = 849.3125
$ws.Range("J39").Value = 2668
$ws.Range("K39").Value = 2547.9375
$ws.Range("L39").Value = 8004
$ws.Range("M39").Value = -2253.9375
$ws.Range("N39").Value = -8592
$ws.Range("H129").Value = 17544724
$ws.Range("I129").Value = 23810114
$ws.Range("J129").Value = 1632
$ws.Range("K129").Value = 71430342
$ws.Range("L129").Value = 4896
$ws.Range("M129").Value = -71425342
$ws.Range("N129").Value = -14896
$ws.Range("H131").Value = 38462908
$ws.Range("J131").Value = 1600.2858
$ws.Range("L131").Value = 4800.857400000001
$ws.Range("N131").Value = -14880.8574

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6874.6
$ws.Range("I102").Value = 1718.25
$ws.Range("K102").Value = 1718.25
$ws.Range("M102").Value = -96.25
$ws.Range("H126").Value = 2426.8572
$ws.Range("I126").Value = 2331.3333
$ws.Range("K126").Value = 6993.999899999999
$ws.Range("M126").Value = -4523.999899999999
$ws.Range("H132").Value = 3013
$ws.Range("I132").Value = 2012
$ws.Range("J132").Value = 4014
$ws.Range("K132").Value = 6036
$ws.Range("L132").Value = 12042
$ws.Range("M132").Value = -3506
$ws.Range("N132").Value = -17102

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1249.5555
$ws.Range("I16").Value = 1299.4667
$ws.Range("K16").Value = 1299.4667
$ws.Range("M16").Value = -1129.4667
$ws.Range("H22").Value = 2949.8572
$ws.Range("J22").Value = 3160
$ws.Range("L22").Value = 3160
$ws.Range("N22").Value = -3750
$ws.Range("H27").Value = 2949.8572
$ws.Range("J27").Value = 3160
$ws.Range("L27").Value = 3160
$ws.Range("N27").Value = -3374
$ws.Range("H61").Value = 3240.1875
$ws.Range("I61").Value = 2703.1428
$ws.Range("K61").Value = 2703.1428
$ws.Range("M61").Value = -2501.1428
$ws.Range("H68").Value = 6038.8
$ws.Range("I68").Value = 5998.2856
$ws.Range("J68").Value = 6133.3335
$ws.Range("K68").Value = 5998.2856
$ws.Range("L68").Value = 6133.3335
$ws.Range("M68").Value = -5249.2856
$ws.Range("N68").Value = -7631.3335
$ws.Range("H71").Value = 6038.8
$ws.Range("I71").Value = 5998.2856
$ws.Range("J71").Value = 6133.3335
$ws.Range("K71").Value = 29991.428
$ws.Range("L71").Value = 30666.6675
$ws.Range("M71").Value = -26247.428
$ws.Range("N71").Value = -38154.6675
$ws.Range("H113").Value = 3240.1875
$ws.Range("I113").Value = 2703.1428
$ws.Range("K113").Value = 2703.1428
$ws.Range("M113").Value = -533.1428000000001
$ws.Range("H132").Value = 1763.08
$ws.Range("I132").Value = 1629.5217
$ws.Range("K132").Value = 4888.5651
$ws.Range("M132").Value = -2358.5651

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 339965.84
$ws.Range("I18").Value = 1000950
$ws.Range("K18").Value = 1000950
$ws.Range("M18").Value = -1000777
$ws.Range("H113").Value = 1156.4375
$ws.Range("J113").Value = 829.25
$ws.Range("L113").Value = 2487.75
$ws.Range("N113").Value = -6827.75
$ws.Range("H132").Value = 2496.5454
$ws.Range("I132").Value = 2367.2903
$ws.Range("K132").Value = 7101.8709
$ws.Range("M132").Value = -4571.8709

Write-Host "Applied Famfrit_Profits market-data refresh."